$wb = $excel.ActiveWorkbook

# ---- Sheet: n1_d40 ----
$ws = $wb.Worksheets.Item("n1_d40")
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Range("F2").Value = 32.3935
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 32.3935
$ws.Range("I2").Value = 5002.2
$ws.Range("F3").Value = 31.9932
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 31.9932
$ws.Range("I3").Value = 5022.5
$ws.Range("F4").Value = 32.2098
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 32.2098
$ws.Range("I4").Value = 5047.9
$ws.Range("F5").Value = 32.0937
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 32.0937
$ws.Range("I5").Value = 5025.4
$ws.Range("F6").Value = 32.1011
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 32.1011
$ws.Range("I6").Value = 5017.9
$ws.Range("F7").Value = 32.1895
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 32.1895
$ws.Range("I7").Value = 5041.6
$ws.Range("F8").Value = 32.9063
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 32.9063
$ws.Range("I8").Value = 5048.6
$ws.Range("F9").Value = 32.4301
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 32.4301
$ws.Range("I9").Value = 5001
$ws.Range("F10").Value = 32.4715
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 32.4715
$ws.Range("I10").Value = 5013.9
$ws.Range("F11").Value = 31.786
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 31.786
$ws.Range("I11").Value = 5023.8
$ws.Range("F12").Value = 32.25747
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 32.25747
$ws.Range("I12").Value = 5024.48

# ---- Sheet: n1_d60 ----
$ws = $wb.Worksheets.Item("n1_d60")
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 5001.9
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 5070.5
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 5045.9
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 5067.8
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 5010.9
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 5078.4
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5055.6
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 5068.9
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 5034.1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 5005.3
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 5043.929999999999

# ---- Sheet: n1_d80 ----
$ws = $wb.Worksheets.Item("n1_d80")
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 5079.8
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 5010.5
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 5080
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 5061.2
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 5046.2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 5018.3
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5047.4
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 5051.8
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 5016.2
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 5058.8
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 5047.02

# ---- Sheet: n1_d100 ----
$ws = $wb.Worksheets.Item("n1_d100")
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(7).Insert()

$ws.Range("F1").Value = "encoding_time_sec"
$ws.Range("G1").Value = "solving_time_sec"
$ws.Range("H1").Value = "total_time_sec"

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 5034.6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 5062.2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 5050.5
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 5049
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 5053.5
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 5028.6
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 5032.6
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 5053.1
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 5072.8
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 5051.4
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 5048.830000000001
